$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Delete row 3 entirely (the "0.03" percentile row); subsequent rows shift up
# and keep referring to the same (unshifted) Sheet1 rows, matching the diff.
$ws.Rows("3:3").Delete()

# Update the active selection to match the new state.
$ws.Activate()
$ws.Range("B7").Select()
